# Edit script: applies the "Finalized W13S3 End lecture" changes.
#
# 1) On slide 12 (sldId 486), shape "TextBox 10" (id 11), italicise every
#    paragraph of the operations list except the first intro line
#    ("And the list of operations below:").
# 2) Refresh the cached "datetimeFigureOut" field text (12/4/2023 ->
#    20/4/2023) across every slide layout, the slide master, and the
#    notes master (this mirrors PowerPoint re-caching the auto date
#    field on save on a later day).

$p = $ppt.ActivePresentation

# --- 1) Italicise the operations list on slide 12 --------------------
$s = $p.Slides.Item(12)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -eq 11) {
        $tr = $sh.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($j = 2; $j -le $paraCount; $j++) {
            $para = $tr.Paragraphs($j, 1)
            $para.Font.Italic = $true
        }
    }
}

# --- 2) Refresh cached date field text --------------------------------
function Update-DatePlaceholder($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shape = $shapes.Item($k)
        $isDate = $false
        try {
            if ($shape.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate) {
            $shape.TextFrame.TextRange.Text = "20/4/2023"
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# Notes master
Update-DatePlaceholder $p.NotesMaster.Shapes

Write-Host "done"
